$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, shifting existing rows 44-49 down to 45-50
$ws.Rows.Item(44).Insert()

# Populate the new row 44 with the data from the diff
$ws.Cells.Item(44, 1).Value = 10
$ws.Cells.Item(44, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(44, 3).Value = "La Araucanía"
$ws.Cells.Item(44, 4).Value = 44722
$ws.Cells.Item(44, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(44, 5).Value = 9
$ws.Cells.Item(44, 6).Value = "Fruta"
$ws.Cells.Item(44, 7).Value = 100107
$ws.Cells.Item(44, 8).Value = "Otros"
$ws.Cells.Item(44, 9).Value = 100107001
$ws.Cells.Item(44, 10).Value = "Caqui"
$ws.Cells.Item(44, 11).Value = "Mankaki"
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 35
$ws.Cells.Item(44, 14).Value = 20000
$ws.Cells.Item(44, 15).Value = 20000
$ws.Cells.Item(44, 16).Value = 20000
$ws.Cells.Item(44, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(44, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(44, 19).Value = 1333
$ws.Cells.Item(44, 20).Value = 15
